# Daily data update - extend "July" sheet with data through 29/7 (rows 28-30)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July")

# --- Add the three new daily rows (28, 29, 30) ---
# Copy the date-column format (style s="1") from the last existing row down
# to the three new rows before populating their values.
$ws.Range("A27").Copy()
$ws.Range("A28:A30").PasteSpecial(-4122)

# Row 28 -> 2024-07-27 (serial 45500)
$ws.Range("A28").Value2 = 45500
$ws.Range("B28").Value2 = 232
$ws.Range("C28").Value2 = 33
$ws.Range("D28").Value2 = 2
$ws.Range("E28").Value2 = 197

# Row 29 -> 2024-07-28 (serial 45501)
$ws.Range("A29").Value2 = 45501
$ws.Range("B29").Value2 = 278
$ws.Range("C29").Value2 = 32
$ws.Range("D29").Value2 = 8
$ws.Range("E29").Value2 = 238

# Row 30 -> 2024-07-29 (serial 45502)
$ws.Range("A30").Value2 = 45502
$ws.Range("B30").Value2 = 277
$ws.Range("C30").Value2 = 35
$ws.Range("D30").Value2 = 5
$ws.Range("E30").Value2 = 237

# --- Extend the summary formulas in I8:K8 to cover the new rows ---
$ws.Range("I8").Formula = "=SUM(B2,B3:B30)"
$ws.Range("J8").Formula = "=SUM(C2,C3:C30)"
$ws.Range("K8").Formula = "=SUM(D2,D3:D30)"

# --- Update the view: scroll back to top and select I8:K8 ---
$ws.Range("A1").Select()
$ws.Range("I8:K8").Select()
